$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.093.43"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.765.02"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.19"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3914"
$ws.Range("E7").Value = "  +2.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3402"
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").Value = "  -3.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07242"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.50"
$ws.Range("E12").Value = "  -3.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9990"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.166"
$ws.Range("E14").Value = "  -4.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.121"
$ws.Range("E15").Value = "  -2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.758.77"
$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001063"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06620"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.47"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9984"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.244"
$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.074.28"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.69"
$ws.Range("E24").Value = "  -3.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.45"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  -3.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.334"
$ws.Range("E28").Value = "  -3.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.952.96"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.286"
$ws.Range("E30").Value = "  -10.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.67"
$ws.Range("E31").Value = "  -3.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.078"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.844"
$ws.Range("E33").Value = "  -4.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08741"
$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.13"
$ws.Range("E35").Value = "  -4.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06223"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02297"
$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.162"
$ws.Range("E38").Value = "  -3.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6528"
$ws.Range("E39").Value = "  -4.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2115"
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("E42").Value = "  -3.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.932"
$ws.Range("E43").Value = "  -3.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9979"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.88"
$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.827"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6027"
$ws.Range("E47").Value = "  -4.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.17"
$ws.Range("E48").Value = "  -4.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.005"
$ws.Range("E49").Value = "  -3.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07021"
$ws.Range("E51").Value = "  -6.42%  "
